$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.040.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.37%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.751.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'602.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'165.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.96%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.750.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.18%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.44%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.66%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'37.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.85%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.380.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.18%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.740.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'69.013.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.48%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'17.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.86%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.81%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +6.79%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'490.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.45%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.725"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.78%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'84.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.42%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -2.85%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.73%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.10%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.86%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.75%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.34%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'31.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.58%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.893.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.35%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.691.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.84%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.140"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +5.42%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.22%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +8.51%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'48.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.61%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.41%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'424.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.94%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'142.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'40.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.42%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +8.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.785.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.02%  "
$ws.Range("E51").Style = "Normal"
